# "Creatinn new branch develop" -- append the next three log rows (47-49)
# to the Hoja1 bootcamp tracker, mirroring what the author typed after
# pasting a block of webex-session data (url/code pairs in D/E) and then
# going back to fill in the C-column topic cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 47 (2024-05-08) ------------------------------------------------
$ws.Range("A47").Value = 45420
$ws.Range("A47").NumberFormat = "d-mmm"
$ws.Range("D47").Value = "https://tecsup.webex.com/tecsup/ldr.php?RCID=0acab4db96997e9e9900eb9c64487d02"
$ws.Range("E47").Value = "8Tu4v6PB"

# --- Row 48 (2024-05-09) ------------------------------------------------
$ws.Range("A48").Value = 45421
$ws.Range("A48").NumberFormat = "d-mmm"
$ws.Range("D48").Value = "https://tecsup.webex.com/tecsup/ldr.php?RCID=f5df756e0325ba56b6b9cf08c928505b"
$ws.Range("E48").Value = "pK5X2sHM"

# --- Row 49 (2024-05-10) -------------------------------------------------
$ws.Range("A49").Value = 45422
$ws.Range("A49").NumberFormat = "d-mmm"
$ws.Range("D49").Value = "https://tecsup.webex.com/tecsup/ldr.php?RCID=5548264b82c24dc9378bd22623e95c9a"
$ws.Range("E49").Value = "iMmwWu8x"

# Topic column, filled in afterwards (matches shared-string append order).
$ws.Range("C47").Value = "Ramas de GIT, "
$ws.Range("C48").Value = "Tutoria"

# Turn the D47 URL into a real hyperlink (new relationship rId29).
$ws.Hyperlinks.Add($ws.Range("D47"), "https://tecsup.webex.com/tecsup/ldr.php?RCID=0acab4db96997e9e9900eb9c64487d02") | Out-Null

# Hyperlinks.Add() stamps a freshly-minted "Hipervinculo" cell format;
# re-apply the workbook's existing hyperlink style (already used by D42..D44)
# via a formats-only paste so no duplicate style entry lingers behind.
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection/scroll position where the author left it.
$ws.Range("C48").Select()
